# Apply updated crypto price/volume data to the worksheet.
# (Generated to reproduce the commit "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds text such as "29.911.48" or "0.7397".
# Force it to stay formatted as Text so Excel does not reinterpret
# values that look like plain numbers (e.g. "0.7397", "1.001") as
# numeric cells when we assign the new strings below.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.911.48'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '1.875.34'
$ws.Range("E3").Value = '  -0.72%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '0.7397'
$ws.Range("E5").Value = '  -3.83%  '
$ws.Range("D6").Value = '242.58'
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '0.3153'
$ws.Range("E8").Value = '  +0.89%  '
$ws.Range("D9").Value = '0.07231'
$ws.Range("E9").Value = '  +0.75%  '
$ws.Range("D10").Value = '24.62'
$ws.Range("E10").Value = '  -4.06%  '
$ws.Range("D11").Value = '0.08333'
$ws.Range("E11").Value = '  -2.76%  '
$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D12").Value = '0.7500'
$ws.Range("E12").Value = '  -1.87%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '5.390'
$ws.Range("E13").Value = '  +0.42%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.866.31'
$ws.Range("E14").Value = '  -3.07%  '
$ws.Range("D15").Value = '92.26'
$ws.Range("E15").Value = '  -1.54%  '
$ws.Range("D16").Value = '29.924.29'
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("D17").Value = '6.099'
$ws.Range("E17").Value = '  -0.92%  '
$ws.Range("D18").Value = '248.06'
$ws.Range("E18").Value = '  +1.34%  '
$ws.Range("D19").Value = '13.56'
$ws.Range("E19").Value = '  -1.62%  '
$ws.Range("D20").Value = '0.000007836'
$ws.Range("E20").Value = '  +0.20%  '
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.18%  '
$ws.Range("D22").Value = '2.142.51'
$ws.Range("E22").Value = '  -2.08%  '
$ws.Range("D23").Value = '8.006'
$ws.Range("E23").Value = '  -0.50%  '
$ws.Range("D24").Value = '1.001'
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").Value = '0.1556'
$ws.Range("E25").Value = '  -5.19%  '
$ws.Range("D26").Value = '9.282'
$ws.Range("E26").Value = '  -1.20%  '
$ws.Range("D27").Value = '165.05'
$ws.Range("E27").Value = '  +1.36%  '
$ws.Range("D28").Value = '18.65'
$ws.Range("E28").Value = '  -0.52%  '
$ws.Range("D29").Value = '2.023'
$ws.Range("E29").Value = '  -0.61%  '
$ws.Range("D30").Value = '1.498'
$ws.Range("E30").Value = '  +2.34%  '
$ws.Range("D31").Value = '4.602'
$ws.Range("E31").Value = '  +1.93%  '
$ws.Range("E32").Value = '  +0.15%  '
$ws.Range("D33").Value = '4.260'
$ws.Range("E33").Value = '  +3.92%  '
$ws.Range("D34").Value = '0.05329'
$ws.Range("E34").Value = '  -2.31%  '
$ws.Range("D35").Value = '1.234'
$ws.Range("E35").Value = '  -0.71%  '
$ws.Range("D36").Value = '0.7505'
$ws.Range("E36").Value = '  +0.72%  '
$ws.Range("D37").Value = '1.000'
$ws.Range("E37").Value = '  -0.18%  '
$ws.Range("D38").Value = '2.698'
$ws.Range("E38").Value = '  -0.10%  '
$ws.Range("D39").Value = '0.01961'
$ws.Range("E39").Value = '  +0.32%  '
$ws.Range("E40").Value = '  -1.06%  '
$ws.Range("D41").Value = '0.4531'
$ws.Range("E41").Value = '  +1.29%  '
$ws.Range("D42").Value = '6.145'
$ws.Range("E42").Value = '  +1.11%  '
$ws.Range("D43").Value = '1.107.41'
$ws.Range("E43").Value = '  +0.09%  '
$ws.Range("D44").Value = '72.27'
$ws.Range("E44").Value = '  -1.32%  '
$ws.Range("D45").Value = '0.8649'
$ws.Range("E45").Value = '  +1.40%  '
$ws.Range("D46").Value = '104.40'
$ws.Range("E46").Value = '  +1.49%  '
$ws.Range("D47").Value = '1.001'
$ws.Range("E47").Value = '  +0.10%  '
$ws.Range("D48").Value = '1.861'
$ws.Range("E48").Value = '  -0.24%  '
$ws.Range("D49").Value = '7.595'
$ws.Range("E49").Value = '  -0.98%  '
$ws.Range("D50").Value = '9.541'
$ws.Range("E50").Value = '  -2.39%  '
$ws.Range("D51").Value = '2.040.15'
$ws.Range("E51").Value = '  -1.24%  '
